$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.214209914207458
$ws.Range("B1").Value = 2.431705236434937
$ws.Range("C1").Value = 7.117269992828369
$ws.Range("D1").Value = 2.256514310836792
$ws.Range("E1").Value = 1.16564416885376
